$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("displaycase_arrowsmith")

# Row 5 (Cheese Curds) had ingredients/allergens still referencing the sandwich item; replace with actual values
$ws.Range("B5").Value = "28% MF Milk / Enzymes / Salt"
$ws.Range("C5").Value = "Milk."

# Replace "waffles" placeholder values in the Nutrition Label column (F2:F5) with "placeholder"
$ws.Range("F2").Value = "placeholder"
$ws.Range("F3").Value = "placeholder"
$ws.Range("F4").Value = "placeholder"
$ws.Range("F5").Value = "placeholder"

# Row 5 no longer needs the extra wrapped height now that the text is shorter
$ws.Rows.Item(5).AutoFit() | Out-Null

# Update the active selection to reflect where the editor last clicked
$ws.Range("F4").Select() | Out-Null
